# Applies commit "Alterada a etapa 12": updates wording of a few test
# steps (client-type selection, quantity, and expected results) across
# the TC1-TC5 test-case blocks on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC1 (rows 10-21 -> steps 1-12)
$ws.Range("D15").Value = "SYSTEM apresenta campos: Tipo de Cliente e Quantidade"
$ws.Range("B16").Value = "Usuário do Sistema altera para tipo de cliente B"
$ws.Range("B17").Value = "Usuário do Sistema informa a quantidade de produtos menor ou igual a zero"
$ws.Range("D17").Value = "SYSTEM exibe mensagem 'A quantidade informada deve ser maior ou igual a 01 (um)!' (MSG002)"

# TC2 (rows 29-40 -> steps 1-12)
$ws.Range("D34").Value = "SYSTEM apresenta campos: Tipo de Cliente e Quantidade"
$ws.Range("B36").Value = "Usuário do Sistema informa a quantidade de produtos menor que 100 unidades"
$ws.Range("D36").Value = "SYSTEM aplica fator de desconto para quantidade < 100: Cliente A (0,90), B (0,85), C (0,80)"

# TC3 (rows 48-59 -> steps 1-12)
$ws.Range("D53").Value = "SYSTEM apresenta campos: Tipo de Cliente e Quantidade"
$ws.Range("B54").Value = "Usuário do Sistema altera para tipo de cliente C"

# TC4 (rows 67-78 -> steps 1-12)
$ws.Range("D72").Value = "SYSTEM apresenta campos: Tipo de Cliente e Quantidade"

# TC5 (rows 86-97 -> steps 1-12)
$ws.Range("D91").Value = "SYSTEM apresenta campos: Tipo de Cliente e Quantidade"
$ws.Range("B92").Value = "Usuário do Sistema altera para tipo de cliente C"
$ws.Range("B93").Value = "Usuário do Sistema informa a quantidade de produtos entre 100 e 999 unidades"
$ws.Range("D93").Value = "SYSTEM aplica fator de desconto para 100 <= quantidade < 1000: Cliente A (0,95), B (0,90), C (0,85)"
